# Auto-generated: updates cryptos list figures (prices + 1h volume deltas)
# Commit: "Updated cryptos list on Tue Jan 23 03:58:15 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'40.023.89"
$ws.Range("E2").Value = "  -2.94%  "

# Row 3
$ws.Range("D3").Value = "'2.338.82"
$ws.Range("E3").Value = "  -4.15%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'310.18"
$ws.Range("E5").Value = "  -2.11%  "

# Row 6
$ws.Range("D6").Value = "'85.24"
$ws.Range("E6").Value = "  -4.55%  "

# Row 7
$ws.Range("D7").Value = "'0.530"
$ws.Range("E7").Value = "  -2.23%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").Value = "'0.483"
$ws.Range("E9").Value = "  -2.66%  "

# Row 10
$ws.Range("D10").Value = "'0.0810"
$ws.Range("E10").Value = "  -2.92%  "

# Row 11
$ws.Range("D11").Value = "'29.96"
$ws.Range("E11").Value = "  -6.99%  "

# Row 12
$ws.Range("E12").Value = "  +0.94%  "

# Row 13
$ws.Range("D13").Value = "'2.699.28"
$ws.Range("E13").Value = "  -3.93%  "

# Row 14
$ws.Range("E14").Value = "  -4.51%  "

# Row 15
$ws.Range("D15").Value = "'14.72"
$ws.Range("E15").Value = "  -5.73%  "

# Row 16
$ws.Range("D16").Value = "'2.360.62"
$ws.Range("E16").Value = "  -3.54%  "

# Row 17
$ws.Range("E17").Value = "  -2.04%  "

# Row 18
$ws.Range("D18").Value = "'39.988.58"
$ws.Range("E18").Value = "  -2.83%  "

# Row 19
$ws.Range("D19").Value = "'0.0₃0904"
$ws.Range("E19").Value = "  -2.23%  "

# Row 20
$ws.Range("D20").Value = "'6.11"
$ws.Range("E20").Value = "  -2.25%  "

# Row 21
$ws.Range("D21").Value = "'68.09"
$ws.Range("E21").Value = "  -5.80%  "

# Row 22
$ws.Range("D22").Value = "'10.70"
$ws.Range("E22").Value = "  -2.98%  "

# Row 23
$ws.Range("D23").Value = "'235.02"
$ws.Range("E23").Value = "  -0.24%  "

# Row 24
$ws.Range("D24").Value = "'2.55"
$ws.Range("E24").Value = "  -5.06%  "

# Row 25
$ws.Range("E25").Value = "  +0.14%  "

# Row 26
$ws.Range("E26").Value = "  -3.43%  "

# Row 27
$ws.Range("D27").Value = "'23.32"
$ws.Range("E27").Value = "  -3.26%  "

# Row 29
$ws.Range("D29").Value = "'9.29"
$ws.Range("E29").Value = "  -2.83%  "

# Row 30
$ws.Range("D30").Value = "'34.65"
$ws.Range("E30").Value = "  -0.79%  "

# Row 31
$ws.Range("D31").Value = "'153.02"
$ws.Range("E31").Value = "  -2.19%  "

# Row 32
$ws.Range("E32").Value = "  -0.12%  "

# Row 33
$ws.Range("D33").Value = "'5.10"
$ws.Range("E33").Value = "  -3.10%  "

# Row 35
$ws.Range("D35").Value = "'0.0719"
$ws.Range("E35").Value = "  -3.68%  "

# Row 36
$ws.Range("E36").Value = "  -0.89%  "

# Row 37
$ws.Range("E37").Value = "  -4.36%  "

# Row 38
$ws.Range("D38").Value = "'0.0987"
$ws.Range("E38").Value = "  -1.93%  "

# Row 39
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'1.73"
$ws.Range("E39").Value = "  -2.70%  "

# Row 40
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "'15.62"
$ws.Range("E40").Value = "  -6.51%  "

# Row 41
$ws.Range("D41").Value = "'3.87"
$ws.Range("E41").Value = "  +0.27%  "

# Row 42
$ws.Range("D42").Value = "'1.954.29"
$ws.Range("E42").Value = "  -1.88%  "

# Row 43
$ws.Range("E43").Value = "  -4.23%  "

# Row 44
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0264"
$ws.Range("E44").Value = "  -4.21%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'17.55"
$ws.Range("E45").Value = "  -4.94%  "

# Row 46
$ws.Range("D46").Value = "'9.39"
$ws.Range("E46").Value = "  -0.42%  "

# Row 47
$ws.Range("E47").Value = "  -6.29%  "

# Row 48
$ws.Range("D48").Value = "'2.559.62"
$ws.Range("E48").Value = "  -4.14%  "

# Row 49
$ws.Range("E49").Value = "  -2.73%  "

# Row 50
$ws.Range("D50").Value = "'70.53"
$ws.Range("E50").Value = "  -3.73%  "

# Row 51
$ws.Range("D51").Value = "'50.80"
$ws.Range("E51").Value = "  -1.77%  "

